{"js": "const replacements = [\n  { oldText: \"26\u00d729=754\", newText: \"96\u00d757=5472\" },\n  { oldText: \"45\u00d711=495\", newText: \"13\u00d786=1118\" },\n  { oldText: \"40\u00d782=3280\", newText: \"56\u00d778=4368\" },\n  { oldText: \"88\u00d777=6776\", newText: \"56\u00d723=1288\" },\n  { oldText: \"54\u00d751=2754\", newText: \"37\u00d773=2701\" },\n  { oldText: \"47\u00d729=1363\", newText: \"33\u00d760=1980\" },\n  { oldText: \"79\u00d781=6399\", newText: \"83\u00d740=3320\" },\n  { oldText: \"83\u00d783=6889\", newText: \"84\u00d727=2268\" },\n  { oldText: \"60\u00d778=4680\", newText: \"77\u00d737=2849\" },\n  { oldText: \"25\u00d771=1775\", newText: \"89\u00d758=5162\" },\n  { oldText: \"78\u00d739=3042\", newText: \"12\u00d788=1056\" },\n  { oldText: \"51\u00d798=4998\", newText: \"39\u00d713=507\" },\n  { oldText: \"46\u00d725=1150\", newText: \"24\u00d724=576\" },\n  { oldText: \"16\u00d768=1088\", newText: \"46\u00d736=1656\" },\n  { oldText: \"87\u00d714=1218\", newText: \"99\u00d738=3762\" },\n  { oldText: \"40\u00d740=1600\", newText: \"25\u00d754=1350\" },\n  { oldText: \"55\u00d750=2750\", newText: \"73\u00d746=3358\" },\n  { oldText: \"19\u00d770=1330\", newText: \"26\u00d764=1664\" },\n  { oldText: \"36\u00d741=1476\", newText: \"87\u00d741=3567\" },\n  { oldText: \"77\u00d746=3542\", newText: \"14\u00d771=994\" },\n  { oldText: \"37\u00d772=2664\", newText: \"22\u00d729=638\" },\n  { oldText: \"23\u00d711=253\", newText: \"74\u00d790=6660\" },\n  { oldText: \"95\u00d788=8360\", newText: \"67\u00d794=6298\" },\n  { oldText: \"77\u00d778=6006\", newText: \"27\u00d729=783\" },\n  { oldText: \"78\u00d766=5148\", newText: \"14\u00d765=910\" },\n];\n\nconst body = context.document.body;\n\nfor (const { oldText, newText } of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"26\u00d729=754\"; New = \"96\u00d757=5472\" }\n    @{ Old = \"45\u00d711=495\"; New = \"13\u00d786=1118\" }\n    @{ Old = \"40\u00d782=3280\"; New = \"56\u00d778=4368\" }\n    @{ Old = \"88\u00d777=6776\"; New = \"56\u00d723=1288\" }\n    @{ Old = \"54\u00d751=2754\"; New = \"37\u00d773=2701\" }\n    @{ Old = \"47\u00d729=1363\"; New = \"33\u00d760=1980\" }\n    @{ Old = \"79\u00d781=6399\"; New = \"83\u00d740=3320\" }\n    @{ Old = \"83\u00d783=6889\"; New = \"84\u00d727=2268\" }\n    @{ Old = \"60\u00d778=4680\"; New = \"77\u00d737=2849\" }\n    @{ Old = \"25\u00d771=1775\"; New = \"89\u00d758=5162\" }\n    @{ Old = \"78\u00d739=3042\"; New = \"12\u00d788=1056\" }\n    @{ Old = \"51\u00d798=4998\"; New = \"39\u00d713=507\" }\n    @{ Old = \"46\u00d725=1150\"; New = \"24\u00d724=576\" }\n    @{ Old = \"16\u00d768=1088\"; New = \"46\u00d736=1656\" }\n    @{ Old = \"87\u00d714=1218\"; New = \"99\u00d738=3762\" }\n    @{ Old = \"40\u00d740=1600\"; New = \"25\u00d754=1350\" }\n    @{ Old = \"55\u00d750=2750\"; New = \"73\u00d746=3358\" }\n    @{ Old = \"19\u00d770=1330\"; New = \"26\u00d764=1664\" }\n    @{ Old = \"36\u00d741=1476\"; New = \"87\u00d741=3567\" }\n    @{ Old = \"77\u00d746=3542\"; New = \"14\u00d771=994\" }\n    @{ Old = \"37\u00d772=2664\"; New = \"22\u00d729=638\" }\n    @{ Old = \"23\u00d711=253\"; New = \"74\u00d790=6660\" }\n    @{ Old = \"95\u00d788=8360\"; New = \"67\u00d794=6298\" }\n    @{ Old = \"77\u00d778=6006\"; New = \"27\u00d729=783\" }\n    @{ Old = \"78\u00d766=5148\"; New = \"14\u00d765=910\" }\n)\n\nforeach ($rep in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($rep.Old, $false, $false, $false, $false, $false, $true, 1, $false, $rep.New, 2)\n}\n"}
